# Apply latest cryptos snapshot values (prices / 1h volume %) scraped from coinranking.com
# Columns: B=Coin name, C=Link, D=Price, E=Volume(1h)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + "30.968.62"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'" + "  +3.55%  "
$ws.Range("E2").ClearFormats()
$ws.Range("D3").Value = "'" + "1.683.29"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'" + "  +3.08%  "
$ws.Range("E3").ClearFormats()
$ws.Range("D4").Value = "'" + "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "'" + "  +0.08%  "
$ws.Range("E4").ClearFormats()
$ws.Range("D5").Value = "'" + "220.12"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'" + "  +2.32%  "
$ws.Range("E5").ClearFormats()
$ws.Range("D6").Value = "'" + "0.531"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "'" + "  +1.94%  "
$ws.Range("E6").ClearFormats()
$ws.Range("D7").Value = "'" + "1.00"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "'" + "  +0.11%  "
$ws.Range("E7").ClearFormats()
$ws.Range("D8").Value = "'" + "29.27"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "'" + "  +2.10%  "
$ws.Range("E8").ClearFormats()
$ws.Range("E9").Value = "'" + "  +2.38%  "
$ws.Range("E9").ClearFormats()
$ws.Range("D10").Value = "'" + "0.0637"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "'" + "  +4.41%  "
$ws.Range("E10").ClearFormats()
$ws.Range("D11").Value = "'" + "0.0908"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "'" + "  +0.78%  "
$ws.Range("E11").ClearFormats()
$ws.Range("D12").Value = "'" + "1.925.14"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "'" + "  +3.16%  "
$ws.Range("E12").ClearFormats()
$ws.Range("D13").Value = "'" + "1.687.82"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "'" + "  +3.32%  "
$ws.Range("E13").ClearFormats()
$ws.Range("D14").Value = "'" + "10.16"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "'" + "  +7.43%  "
$ws.Range("E14").ClearFormats()
$ws.Range("E15").Value = "'" + "  +5.04%  "
$ws.Range("E15").ClearFormats()
$ws.Range("D16").Value = "'" + "4.12"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "'" + "  +7.21%  "
$ws.Range("E16").ClearFormats()
$ws.Range("D17").Value = "'" + "30.978.84"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "'" + "  +3.58%  "
$ws.Range("E17").ClearFormats()
$ws.Range("D18").Value = "'" + "66.69"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "'" + "  +1.85%  "
$ws.Range("E18").ClearFormats()
$ws.Range("D19").Value = "'" + "246.85"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "'" + "  +2.37%  "
$ws.Range("E19").ClearFormats()
$ws.Range("D20").Value = "'" + "0.0" + [char]0x2083 + "0720"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "'" + "  +2.31%  "
$ws.Range("E20").ClearFormats()
$ws.Range("D21").Value = "'" + "1.00"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "'" + "  +0.12%  "
$ws.Range("E21").ClearFormats()
$ws.Range("D22").Value = "'" + "4.28"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "'" + "  +3.47%  "
$ws.Range("E22").ClearFormats()
$ws.Range("E23").Value = "'" + "  +1.62%  "
$ws.Range("E23").ClearFormats()
$ws.Range("E24").Value = "'" + "  -0.98%  "
$ws.Range("E24").ClearFormats()
$ws.Range("D25").Value = "'" + "158.68"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "'" + "  +0.62%  "
$ws.Range("E25").ClearFormats()
$ws.Range("E26").Value = "'" + "  +2.64%  "
$ws.Range("E26").ClearFormats()
$ws.Range("D27").Value = "'" + "0.112"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "'" + "  +2.38%  "
$ws.Range("E27").ClearFormats()
$ws.Range("D28").Value = "'" + "6.69"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "'" + "  +0.96%  "
$ws.Range("E28").ClearFormats()
$ws.Range("D29").Value = "'" + "1.00"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "'" + "  +0.15%  "
$ws.Range("E29").ClearFormats()
$ws.Range("E30").Value = "'" + "  +2.22%  "
$ws.Range("E30").ClearFormats()
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "'" + "1.15"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "'" + "  +3.67%  "
$ws.Range("E31").ClearFormats()
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "'" + "3.51"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "'" + "  +3.63%  "
$ws.Range("E32").ClearFormats()
$ws.Range("D33").Value = "'" + "3.33"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "'" + "  +4.64%  "
$ws.Range("E33").ClearFormats()
$ws.Range("D34").Value = "'" + "1.515.32"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "'" + "  +6.26%  "
$ws.Range("E34").ClearFormats()
$ws.Range("E35").Value = "'" + "  +2.73%  "
$ws.Range("E35").ClearFormats()
$ws.Range("E36").Value = "'" + "  +0.71%  "
$ws.Range("E36").ClearFormats()
$ws.Range("D37").Value = "'" + "83.27"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "'" + "  +10.26%  "
$ws.Range("E37").ClearFormats()
$ws.Range("D38").Value = "'" + "0.612"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "'" + "  +9.95%  "
$ws.Range("E38").ClearFormats()
$ws.Range("D39").Value = "'" + "0.0179"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "'" + "  +4.83%  "
$ws.Range("E39").ClearFormats()
$ws.Range("D40").Value = "'" + "2.70"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "'" + "  -3.11%  "
$ws.Range("E40").ClearFormats()
$ws.Range("D41").Value = "'" + "2.30"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "'" + "  +0.48%  "
$ws.Range("E41").ClearFormats()
$ws.Range("D42").Value = "'" + "2.03"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "'" + "  +2.03%  "
$ws.Range("E42").ClearFormats()
$ws.Range("D43").Value = "'" + "0.841"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "'" + "  +0.56%  "
$ws.Range("E43").ClearFormats()
$ws.Range("D44").Value = "'" + "0.0502"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "'" + "  +0.22%  "
$ws.Range("E44").ClearFormats()
$ws.Range("E45").Value = "'" + "  +2.52%  "
$ws.Range("E45").ClearFormats()
$ws.Range("D46").Value = "'" + "1.00"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "'" + "  +0.14%  "
$ws.Range("E46").ClearFormats()
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").Value = "'" + "5.57"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "'" + "  +4.33%  "
$ws.Range("E47").ClearFormats()
$ws.Range("B48").Value = "BitcoinSV"
$ws.Range("C48").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D48").Value = "'" + "51.74"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "'" + "  +6.85%  "
$ws.Range("E48").ClearFormats()
$ws.Range("D49").Value = "'" + "1.819.50"
$ws.Range("D49").ClearFormats()
$ws.Range("D50").Value = "'" + "0.0" + [char]0x2086 + "0117"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "'" + "  +6.13%  "
$ws.Range("E50").ClearFormats()
$ws.Range("D51").Value = "'" + "93.51"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "'" + "  +1.12%  "
$ws.Range("E51").ClearFormats()
